$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The workbook tracks 8 fiscal-year columns (D:K) of financial data for the
# ticker, per section (Income Statement / Balance Sheet / Cash Flow
# Statement). This update adds one more (newer) fiscal year of data, which
# pushes the existing data one column to the right (old D:K -> new E:L) and
# populates the freed-up column D with the new year's figures.
# ---------------------------------------------------------------------------

# 1. Insert a new blank column at D - this shifts existing D:K data to E:L,
#    carrying over each cell's formatting/style automatically.
$ws.Columns("D:D").Insert()

# 2. The freshly inserted column D is unformatted; copy number formats
#    (date format for the "Period Ending" rows, #,##0 for the numeric rows)
#    from column E (which now holds what used to be column D) so the new
#    column matches the rest of the table. Rows 5, 6, 37 and 79 are
#    single-cell section headers (no D:K data there), so each data block
#    is handled with its own Copy/PasteSpecial to avoid touching them.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# 3. Populate column D with the new fiscal year's figures.

# --- Income Statement ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 5758000
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -189000
$ws.Range("D17").Value = 1552000
$ws.Range("D18").Value = 4206000
$ws.Range("D20").Value = -2023000
$ws.Range("D21").Value = 2672000
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 2183000
$ws.Range("D24").Value = 491000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1692000
$ws.Range("D27").Value = 1663000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 29000
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 2023000
$ws.Range("D33").Value = 1692000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1692000

# --- Balance Sheet ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 4074000
$ws.Range("D42").Value = 1480000
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 791000
$ws.Range("D49").Value = 7775000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 160518000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 14433000
$ws.Range("D62").Value = 573000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 139701000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 247000
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 5385000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 19977000
$ws.Range("D77").Value = 0

# --- Cash Flow Statement ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 1692000
$ws.Range("D83").Value = 489000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1767000
$ws.Range("D91").Value = -232000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -7077000
$ws.Range("D96").Value = -485000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 6352000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 1042000
